$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'309.15"
$ws.Range("E2").Value = "'-3.94%"
$ws.Range("D3").Value = "'40.10"
$ws.Range("E3").Value = "'-6.19%"
$ws.Range("E4").Value = "'-0.37%"
$ws.Range("D5").Value = "'0.07746"
$ws.Range("E5").Value = "'-5.41%"
$ws.Range("D6").Value = "'4.259"
$ws.Range("E6").Value = "'-0.54%"
$ws.Range("D7").Value = "'1.630"
$ws.Range("E7").Value = "'-9.40%"
$ws.Range("D8").Value = "'0.8820"
$ws.Range("E8").Value = "'-5.33%"
$ws.Range("D9").Value = "'0.1026"
$ws.Range("E9").Value = "'-7.51%"
$ws.Range("D10").Value = "'0.1748"
$ws.Range("E10").Value = "'-6.20%"
$ws.Range("D11").Value = "'0.09075"
$ws.Range("E11").Value = "'-4.02%"
$ws.Range("D12").Value = "'0.04443"
$ws.Range("E12").Value = "'-4.14%"
$ws.Range("E13").Value = "'-0.25%"
$ws.Range("D14").Value = "'0.001257"
$ws.Range("E14").Value = "'-2.83%"
$ws.Range("D15").Value = "'0.005818"
$ws.Range("E15").Value = "'0.97%"
$ws.Range("E16").Value = "'2,413.19%"
$ws.Range("D17").Value = "'3.356"
$ws.Range("E17").Value = "'-0.25%"
$ws.Range("D18").Value = "'2.419"
$ws.Range("E18").Value = "'-4.31%"
$ws.Range("E19").Value = "'-2.92%"
$ws.Range("D20").Value = "'7.038"
$ws.Range("E20").Value = "'-4.85%"
$ws.Range("D21").Value = "'0.1340"
$ws.Range("E21").Value = "'-3.55%"
$ws.Range("D22").Value = "'0.2788"
$ws.Range("E22").Value = "'10.57%"
$ws.Range("D23").Value = "'0.04177"
$ws.Range("E23").Value = "'0.39%"
$ws.Range("E24").Value = "'-3.50%"
$ws.Range("D25").Value = "'0.004078"
$ws.Range("E25").Value = "'-8.41%"
$ws.Range("E26").Value = "'8.36%"
$ws.Range("D38").Value = "'0.02373"
$ws.Range("E38").Value = "'-13.65%"
$ws.Range("D39").Value = "'0.05220"
$ws.Range("E39").Value = "'-6.49%"
$ws.Range("D40").Value = "'0.007948"
$ws.Range("E40").Value = "'-0.79%"
$ws.Range("E41").Value = "'-4.99%"
$ws.Range("D42").Value = "'0.006354"
$ws.Range("E42").Value = "'-2.94%"
$ws.Range("D43").Value = "'0.001966"
$ws.Range("E43").Value = "'-5.79%"
$ws.Range("D44").Value = "'0.008759"
$ws.Range("E44").Value = "'15.90%"
$ws.Range("D45").Value = "'0.3345"
$ws.Range("E45").Value = "'-4.31%"
$ws.Range("D46").Value = "'0.00006558"
$ws.Range("E46").Value = "'-5.90%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.03%"
$ws.Range("E48").Value = "'98.35%"
$ws.Range("D49").Value = "'0.004481"
$ws.Range("E49").Value = "'28.85%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'0.03%"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'0.03%"
